# Edit script: apply "Se agregó folio, validación y limpieza de datos" changes
# to "Base General.xlsx"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($range, [string]$value)
    # Force the cell to keep its content as literal text instead of letting
    # Excel auto-coerce numeric-looking strings (phone numbers, postal codes,
    # amounts prefixed with "$") into Number/Currency values.
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

# --- Header row updates -------------------------------------------------
# G1: "Turno" -> "Hora" ; H1: "Hora" -> "Folio" (new Folio column)
$ws.Range("G1").Value = "Hora"
$ws.Range("H1").Value = "Folio"

# --- Row 2: data cleanup / uppercase / folio ----------------------------
$ws.Range("A2").Value = "ALEX SERRANO DURÁN"
Set-TextValue $ws.Range("B2") "5563193656"
Set-TextValue $ws.Range("C2") "alexserrano0805@gmail.com"
Set-TextValue $ws.Range("D2") "52950"
$ws.Range("E2").Value = "CONSULTA"
Set-TextValue $ws.Range("F2") "`$600"
$ws.Range("G2").Value = "04:24"
$ws.Range("H2").Value = "000-20220825M"
$ws.Range("I2").Value = "26/08/2022"

# --- Row 3: data cleanup / uppercase / folio ----------------------------
$ws.Range("A3").Value = "ALEX SERRANO DURÁN"
Set-TextValue $ws.Range("B3") "5563193656"
Set-TextValue $ws.Range("C3") "alexserrano0805@gmail.com"
Set-TextValue $ws.Range("D3") "52950"
$ws.Range("E3").Value = "CONSULTA"
Set-TextValue $ws.Range("F3") "`$400"
$ws.Range("G3").Value = "04:26"
$ws.Range("H3").Value = "001-20220825N"
$ws.Range("I3").Value = "26/08/2022"

# --- Remove stale/incomplete rows 4-8 (data cleanup) --------------------
$ws.Rows("4:8").Delete()

# --- Column width tweaks -------------------------------------------------
$ws.Columns.Item(2).ColumnWidth = 10.998697916666666
$ws.Columns.Item(3).ColumnWidth = 11.666666666666666
$ws.Columns.Item(4).ColumnWidth = 39.666666666666664
$ws.Columns.Item(5).ColumnWidth = 26.998697916666668
$ws.Columns.Item(6).ColumnWidth = 10.330729166666666

# --- Selection moves to H2 (new Folio column) ----------------------------
$ws.Range("H2").Select() | Out-Null
